{"js": "const body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst items = body.paragraphs.items;\nconst firstPara = items[0];\nconst lastPara = items[items.length - 1];\nconst targetRange = firstPara.getRange(\"Start\").expandTo(lastPara.getRange(\"End\"));\n\nconst newBodyFragment = `<w:p w14:paraId=\"0F5C251B\" w14:textId=\"3E94BEED\" w:rsidR=\"00F054DE\" w:rsidRDefault=\"00520097\"><w:r><w:t>I started on the Second Brief called Speedometer</w:t></w:r><w:r><w:t>. It</w:t></w:r><w:r><w:t xml:space=\"preserve\"> allows you to know how fast the car is driving by showing the Kilometres by hour on the left side of the screen with a minimum of 0 and a maximum of \u201c260\u201d</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> but it goes over that range when you fall off the pathway. I </w:t></w:r><w:r><w:t xml:space=\"preserve\">followed a tutorial </w:t></w:r><w:r><w:t>to</w:t></w:r><w:r><w:t xml:space=\"preserve\"> have a Rigidbody to </w:t></w:r><w:r><w:t>determine</w:t></w:r><w:r><w:t xml:space=\"preserve\"> if the speedometer works</w:t></w:r><w:r><w:t xml:space=\"preserve\">. This </w:t></w:r><w:r><w:t xml:space=\"preserve\">is </w:t></w:r><w:r><w:t>done by the scripts below.</w:t></w:r></w:p><w:p w14:paraId=\"6B6D6D13\" w14:textId=\"47AC2661\" w:rsidR=\"00520097\" w:rsidRDefault=\"00520097\"><w:r><w:t xml:space=\"preserve\">Create a </w:t></w:r><w:r><w:t>script allowing the player to move down the road and interact with things, such as</w:t></w:r><w:r><w:t xml:space=\"preserve\"> obstacles. </w:t></w:r><w:r><w:t>It is called \u201c</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>PlayerController</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>\u201d, allowing you to move forward horizontally and vertically down the road. It</w:t></w:r><w:r><w:t xml:space=\"preserve\"> also </w:t></w:r><w:r><w:t>lets the game</w:t></w:r><w:r><w:t xml:space=\"preserve\"> know the player has </w:t></w:r><w:r><w:t>advanced,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> so the speedometer should add force to </w:t></w:r><w:r><w:t xml:space=\"preserve\">the </w:t></w:r><w:r><w:t>vehicle</w:t></w:r><w:r><w:t>. Hence,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the speedometer knows you\u2019re moving. </w:t></w:r></w:p><w:p w14:paraId=\"71CB9133\" w14:textId=\"65F0C18A\" w:rsidR=\"00DC733E\" w:rsidRDefault=\"00DC733E\"><w:r><w:t>The second Script will be the \u201c</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>FollowPlayer</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>\u201d</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> which allows the Main Camera to follow the play a certain distance away from the Vehicle so you can know what </w:t></w:r><w:r><w:t>is happening</w:t></w:r><w:r><w:t xml:space=\"preserve\"> to the car you are always driving.</w:t></w:r></w:p><w:p w14:paraId=\"2B02DD2D\" w14:textId=\"5838D6E3\" w:rsidR=\"00DC733E\" w:rsidRDefault=\"00DC733E\"><w:r><w:t>The last Script will be the \u201cSpeedometer\u201d Scripts that allows me to</w:t></w:r><w:r><w:t xml:space=\"preserve\"> make the Speedometer on the left side work by giving a rigidbody to latch on to and a minimum and maximum speed so the character can\u2019t go beyond it and go too fast</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> I also add the text and image into the script</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> so they follow along with the car when </w:t></w:r><w:r><w:t>it\u2019s</w:t></w:r><w:r><w:t xml:space=\"preserve\"> driving.</w:t></w:r><w:r><w:t xml:space=\"preserve\">  </w:t></w:r></w:p><w:p><w:r><w:t>To make a speedometer happen, you must first have a working vehicle and its Rigidbody working, then a picture of a Speedometer and an arrow that goes up and down when speeding with the car</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>which then you go on the hierarchy and right click and find UI and click the raw image that allows you paste your pictures then you adjust however you like. After that, you make a speedometer script, and in there, go to the very top where it says using system or engine</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> go on a new line and write (using Unity</w:t></w:r><w:r><w:t>Engine.UI) to tell the game that you are referring to it in your project. Then</w:t></w:r><w:r><w:t xml:space=\"preserve\">, in the public class with the parentheses at the bottom, press </w:t></w:r><w:r><w:t>enters</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">and </w:t></w:r><w:r><w:t xml:space=\"preserve\">write </w:t></w:r><w:r><w:t xml:space=\"preserve\">public </w:t></w:r><w:r><w:t xml:space=\"preserve\">Rigidbody target (Car), </w:t></w:r><w:r><w:t xml:space=\"preserve\">to let the engine know that the speed you want to track is the car. </w:t></w:r></w:p><w:p><w:r><w:t>Adds</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the public float</w:t></w:r><w:r><w:t xml:space=\"preserve\"> and calls it max speed, which will be </w:t></w:r><w:r><w:t>zero,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>so you see the speed go up the on your speedometer</w:t></w:r><w:r><w:t xml:space=\"preserve\">. Then you add </w:t></w:r><w:r><w:t>the two lines underneath, which will also be a public float called min speed arrow angle and Max speed arrow angle so that the car doesn\u2019t exceed the speed limit on the speedometer and doesn\u2019t go below your speedometer either. You add a [</w:t></w:r><w:r><w:t>Header (</w:t></w:r><w:r><w:t>UI</w:t></w:r><w:r><w:t>)] just</w:t></w:r><w:r><w:t xml:space=\"preserve\"> like this so you can refer to UI in the game, and below that, you add a public float that can refer you to the text where you have the speed label that tells you how many kilometres you\u2019re going at</w:t></w:r><w:r><w:t>. Public RectTransform arrow so it can refer to the hand in the UI, an indicator that follows along the numbers in the speedometer, so you see how your speed is in the game. You then add a private float to the speed game with the numbers and f to let it know that you are using a float. Then in the private void update, which might need to be typed in.  In the private void (speed = target.velocity.magnitude * 3.6f;)</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> which will convert that into kilometres</w:t></w:r><w:r><w:t xml:space=\"preserve\">, tell </w:t></w:r><w:r><w:t>you about</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>the momentum in the Rigidbody and control how the car goes.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before=\"240\"/></w:pPr><w:r><w:t xml:space=\"preserve\">Add an if(){}, which looks like this. </w:t></w:r><w:r><w:t>Which</w:t></w:r><w:r><w:t xml:space=\"preserve\"> means if I do this, then do that. In the first bracket, you add speed Label !=  null</w:t></w:r><w:r><w:t>; you do not refer to any object, in the parentheses. You add a speed Label. text((int) (integer) + \u201ckm/h; which means you refer to the text in the speed label image and another if() inside the brackets</w:t></w:r><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:r><w:t xml:space=\"preserve\">you add an (arrow != null) underneath that you add </w:t></w:r><w:r><w:t>arrow</w:t></w:r><w:r><w:t xml:space=\"preserve\">.localEulerAngles </w:t></w:r><w:r><w:t xml:space=\"preserve\">(allows you coordinate with the parents object rotations) = new vector(0, 0, mathf.Lerp(minspeedArrowAngle, </w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>ax</w:t></w:r><w:r><w:t>speedArrowAngle</w:t></w:r><w:r><w:t>, speed/ maxSpeed)</w:t></w:r><w:r><w:softHyphen/><w:t xml:space=\"preserve\">) allows you to move or change the values over some time.  If your game doesn\u2019t show the </w:t></w:r><w:r><w:t>arrow\u2019s movement</w:t></w:r><w:r><w:t xml:space=\"preserve\"> on the speedometer</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> go to </w:t></w:r><w:r><w:t xml:space=\"preserve\">the </w:t></w:r><w:r><w:t>player controller</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> a</w:t></w:r><w:r><w:t>dd</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">rb. add force in front of the c# that allows your player to move forward. </w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/>`;\n\nconst flatOpc = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n<w:body>${newBodyFragment}</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntargetRange.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$firstPara = $d.Paragraphs.Item(1)\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$targetRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)\n\n$newBodyFragment = @'\n<w:p w14:paraId=\"0F5C251B\" w14:textId=\"3E94BEED\" w:rsidR=\"00F054DE\" w:rsidRDefault=\"00520097\"><w:r><w:t>I started on the Second Brief called Speedometer</w:t></w:r><w:r><w:t>. It</w:t></w:r><w:r><w:t xml:space=\"preserve\"> allows you to know how fast the car is driving by showing the Kilometres by hour on the left side of the screen with a minimum of 0 and a maximum of \u201c260\u201d</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> but it goes over that range when you fall off the pathway. I </w:t></w:r><w:r><w:t xml:space=\"preserve\">followed a tutorial </w:t></w:r><w:r><w:t>to</w:t></w:r><w:r><w:t xml:space=\"preserve\"> have a Rigidbody to </w:t></w:r><w:r><w:t>determine</w:t></w:r><w:r><w:t xml:space=\"preserve\"> if the speedometer works</w:t></w:r><w:r><w:t xml:space=\"preserve\">. This </w:t></w:r><w:r><w:t xml:space=\"preserve\">is </w:t></w:r><w:r><w:t>done by the scripts below.</w:t></w:r></w:p><w:p w14:paraId=\"6B6D6D13\" w14:textId=\"47AC2661\" w:rsidR=\"00520097\" w:rsidRDefault=\"00520097\"><w:r><w:t xml:space=\"preserve\">Create a </w:t></w:r><w:r><w:t>script allowing the player to move down the road and interact with things, such as</w:t></w:r><w:r><w:t xml:space=\"preserve\"> obstacles. </w:t></w:r><w:r><w:t>It is called \u201c</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>PlayerController</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>\u201d, allowing you to move forward horizontally and vertically down the road. It</w:t></w:r><w:r><w:t xml:space=\"preserve\"> also </w:t></w:r><w:r><w:t>lets the game</w:t></w:r><w:r><w:t xml:space=\"preserve\"> know the player has </w:t></w:r><w:r><w:t>advanced,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> so the speedometer should add force to </w:t></w:r><w:r><w:t xml:space=\"preserve\">the </w:t></w:r><w:r><w:t>vehicle</w:t></w:r><w:r><w:t>. Hence,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the speedometer knows you\u2019re moving. </w:t></w:r></w:p><w:p w14:paraId=\"71CB9133\" w14:textId=\"65F0C18A\" w:rsidR=\"00DC733E\" w:rsidRDefault=\"00DC733E\"><w:r><w:t>The second Script will be the \u201c</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>FollowPlayer</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>\u201d</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> which allows the Main Camera to follow the play a certain distance away from the Vehicle so you can know what </w:t></w:r><w:r><w:t>is happening</w:t></w:r><w:r><w:t xml:space=\"preserve\"> to the car you are always driving.</w:t></w:r></w:p><w:p w14:paraId=\"2B02DD2D\" w14:textId=\"5838D6E3\" w:rsidR=\"00DC733E\" w:rsidRDefault=\"00DC733E\"><w:r><w:t>The last Script will be the \u201cSpeedometer\u201d Scripts that allows me to</w:t></w:r><w:r><w:t xml:space=\"preserve\"> make the Speedometer on the left side work by giving a rigidbody to latch on to and a minimum and maximum speed so the character can\u2019t go beyond it and go too fast</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> I also add the text and image into the script</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> so they follow along with the car when </w:t></w:r><w:r><w:t>it\u2019s</w:t></w:r><w:r><w:t xml:space=\"preserve\"> driving.</w:t></w:r><w:r><w:t xml:space=\"preserve\">  </w:t></w:r></w:p><w:p><w:r><w:t>To make a speedometer happen, you must first have a working vehicle and its Rigidbody working, then a picture of a Speedometer and an arrow that goes up and down when speeding with the car</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>which then you go on the hierarchy and right click and find UI and click the raw image that allows you paste your pictures then you adjust however you like. After that, you make a speedometer script, and in there, go to the very top where it says using system or engine</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> go on a new line and write (using Unity</w:t></w:r><w:r><w:t>Engine.UI) to tell the game that you are referring to it in your project. Then</w:t></w:r><w:r><w:t xml:space=\"preserve\">, in the public class with the parentheses at the bottom, press </w:t></w:r><w:r><w:t>enters</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">and </w:t></w:r><w:r><w:t xml:space=\"preserve\">write </w:t></w:r><w:r><w:t xml:space=\"preserve\">public </w:t></w:r><w:r><w:t xml:space=\"preserve\">Rigidbody target (Car), </w:t></w:r><w:r><w:t xml:space=\"preserve\">to let the engine know that the speed you want to track is the car. </w:t></w:r></w:p><w:p><w:r><w:t>Adds</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the public float</w:t></w:r><w:r><w:t xml:space=\"preserve\"> and calls it max speed, which will be </w:t></w:r><w:r><w:t>zero,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>so you see the speed go up the on your speedometer</w:t></w:r><w:r><w:t xml:space=\"preserve\">. Then you add </w:t></w:r><w:r><w:t>the two lines underneath, which will also be a public float called min speed arrow angle and Max speed arrow angle so that the car doesn\u2019t exceed the speed limit on the speedometer and doesn\u2019t go below your speedometer either. You add a [</w:t></w:r><w:r><w:t>Header (</w:t></w:r><w:r><w:t>UI</w:t></w:r><w:r><w:t>)] just</w:t></w:r><w:r><w:t xml:space=\"preserve\"> like this so you can refer to UI in the game, and below that, you add a public float that can refer you to the text where you have the speed label that tells you how many kilometres you\u2019re going at</w:t></w:r><w:r><w:t>. Public RectTransform arrow so it can refer to the hand in the UI, an indicator that follows along the numbers in the speedometer, so you see how your speed is in the game. You then add a private float to the speed game with the numbers and f to let it know that you are using a float. Then in the private void update, which might need to be typed in.  In the private void (speed = target.velocity.magnitude * 3.6f;)</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> which will convert that into kilometres</w:t></w:r><w:r><w:t xml:space=\"preserve\">, tell </w:t></w:r><w:r><w:t>you about</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>the momentum in the Rigidbody and control how the car goes.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before=\"240\"/></w:pPr><w:r><w:t xml:space=\"preserve\">Add an if(){}, which looks like this. </w:t></w:r><w:r><w:t>Which</w:t></w:r><w:r><w:t xml:space=\"preserve\"> means if I do this, then do that. In the first bracket, you add speed Label !=  null</w:t></w:r><w:r><w:t>; you do not refer to any object, in the parentheses. You add a speed Label. text((int) (integer) + \u201ckm/h; which means you refer to the text in the speed label image and another if() inside the brackets</w:t></w:r><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:r><w:t xml:space=\"preserve\">you add an (arrow != null) underneath that you add </w:t></w:r><w:r><w:t>arrow</w:t></w:r><w:r><w:t xml:space=\"preserve\">.localEulerAngles </w:t></w:r><w:r><w:t xml:space=\"preserve\">(allows you coordinate with the parents object rotations) = new vector(0, 0, mathf.Lerp(minspeedArrowAngle, </w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>ax</w:t></w:r><w:r><w:t>speedArrowAngle</w:t></w:r><w:r><w:t>, speed/ maxSpeed)</w:t></w:r><w:r><w:softHyphen/><w:t xml:space=\"preserve\">) allows you to move or change the values over some time.  If your game doesn\u2019t show the </w:t></w:r><w:r><w:t>arrow\u2019s movement</w:t></w:r><w:r><w:t xml:space=\"preserve\"> on the speedometer</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> go to </w:t></w:r><w:r><w:t xml:space=\"preserve\">the </w:t></w:r><w:r><w:t>player controller</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> a</w:t></w:r><w:r><w:t>dd</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">rb. add force in front of the c# that allows your player to move forward. </w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/>\n'@\n\n$flatOpc = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' + '<w:body>' + $newBodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$targetRange.InsertXML($flatOpc)\n"}
